$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 100, pushing existing rows 100-108 down to 101-109
$ws.Rows.Item(100).Insert()

# Populate the newly inserted row 100 with the new weekly record
$ws.Cells.Item(100, 1).Value = 6
$ws.Cells.Item(100, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(100, 3).Value = "Metropolitana"
$ws.Cells.Item(100, 4).Value = 45154
$ws.Cells.Item(100, 5).Value = 13
$ws.Cells.Item(100, 6).Value = 100112035
$ws.Cells.Item(100, 7).Value = "Bruselas (repollito)"
$ws.Cells.Item(100, 8).Value = "Sin especificar"
$ws.Cells.Item(100, 9).Value = "Primera"
$ws.Cells.Item(100, 10).Value = 400
$ws.Cells.Item(100, 11).Value = 16000
$ws.Cells.Item(100, 12).Value = 18000
$ws.Cells.Item(100, 13).Value = 17150
$ws.Cells.Item(100, 14).Value = "`$/malla 15 kilos"
$ws.Cells.Item(100, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(100, 16).Value = 1143
$ws.Cells.Item(100, 17).Value = 15
$ws.Cells.Item(100, 18).Value = "Hortaliza"
